$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only touch the brand-new cells introduced by this edit (two extra
# translation columns: Onkelos / Jonathan). Every pre-existing cell is
# left untouched so its shared-string reference / style is preserved
# exactly as-is; the sheet dimension/row spans grow automatically.
$ws.Range("J1").Value = "Onkelos"
$ws.Range("K1").Value = "Jonathan"
$ws.Range("J2").Value = "“Go, gather the elders of Yisrael, and say to them, ‘Adonoy, the God of your fathers appeared [<b>became revealed</b>] to me—the God of Avraham, Yitzchok and Yaakov—saying, “I have indeed been mindful of you, regarding that which is being done to you in Egypt."
$ws.Range("K2").Value = "Go, and assemble the elders of Israel, and say to them, The Lord God of your fathers hath appeared unto me, the God of Abraham, Izhak, and Jakob, saying, Remembering, I have remembered you, and the injury that is done you in Mizraim;"

# Match the bold/bordered/centered header formatting already used by
# A1:I1 (copy format from I1, the last existing header cell) so J1/K1
# line up visually with the rest of row 1, reusing the existing style
# (no new cellXfs entries).
$ws.Range("I1").Copy()
$ws.Range("J1:K1").PasteSpecial(-4122)

